$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.078.45'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '1.854.02'
$ws.Range("E3").Value = '  +2.60%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.03%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.90'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.331'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0695'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("D12").Value = '2.121.18'
$ws.Range("E12").Value = '  +2.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = '1.857.32'
$ws.Range("E14").Value = '  +2.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.681'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.82%  '
$ws.Range("E16").Value = '  +2.54%  '
$ws.Range("D17").Value = '35.075.54'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("D19").Value = '0.0₃0797'
$ws.Range("E19").Value = '  +1.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '241.43'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.86%  '
$ws.Range("E22").Value = '  +2.43%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  +1.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +25.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.09%  '
$ws.Range("E29").Value = '  +3.13%  '
$ws.Range("E30").Value = '  +2.37%  '
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("E33").Value = '  +1.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.07'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +15.41%  '
$ws.Range("E35").Value = '  +22.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.787'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +13.03%  '
$ws.Range("E37").Value = '  -2.54%  '
$ws.Range("E38").Value = '  +12.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '91.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0202'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.73%  '
$ws.Range("D41").Value = '1.350.52'
$ws.Range("E41").Value = '  +1.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +71.33%  '
$ws.Range("E45").Value = '  -2.84%  '
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0539'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.89%  '
$ws.Range("D49").Value = '2.034.81'
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +17.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0680'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.44%  '
